$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9750741720199585
$ws.Range("B1").Value = 2.077753067016602
$ws.Range("C1").Value = 7.518755912780762
$ws.Range("D1").Value = 2.511236429214478
$ws.Range("E1").Value = 1.376458644866943
